# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff.
# Updates the Price (D) and Volume(1h) (E) columns for most rows, and
# swaps two pairs of Coin/Link rows (41<->42, 50<->51) to reflect the
# refreshed coinranking.com snapshot used for this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '73.448.26'
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '4.058.76'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '573.86'
$ws.Range('E5').Value = '  +7.83%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '151.89'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '4.053.25'
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.696'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.766'
$ws.Range('E10').Value = '  +2.65%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.28'
$ws.Range('E12').Value = '  +14.09%  '
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.27'
$ws.Range('E14').Value = '  +6.17%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.714.94'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.066.43'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.43'
$ws.Range('E17').Value = '  +3.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.84'
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('E19').Value = '  +3.57%  '
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '73.379.73'
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '446.62'
$ws.Range('E22').Value = '  +4.75%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.60'
$ws.Range('E23').Value = '  +9.77%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '98.63'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.62'
$ws.Range('E25').Value = '  +3.75%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '14.74'
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '4.29'
$ws.Range('E27').Value = '  +20.79%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.52'
$ws.Range('E28').Value = '  +3.94%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '11.13'
$ws.Range('E29').Value = '  +4.71%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.97'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '37.32'
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.88'
$ws.Range('E32').Value = '  +12.85%  '
$ws.Range('E33').Value = '  +4.38%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '13.66'
$ws.Range('E34').Value = '  +2.53%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '688.52'
$ws.Range('E35').Value = '  +2.36%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '48.67'
$ws.Range('E36').Value = '  +14.52%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '68.02'
$ws.Range('E37').Value = '  +3.77%  '
$ws.Range('D38').Value = '0.0₃0916'
$ws.Range('E38').Value = '  +11.01%  '
$ws.Range('E39').Value = '  +5.19%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.151'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.40'
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.26'
$ws.Range('E42').Value = '  +16.86%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0497'
$ws.Range('E45').Value = '  +2.41%  '
$ws.Range('E46').Value = '  +0.24%  '
$ws.Range('E47').Value = '  +1.65%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.75'
$ws.Range('E48').Value = '  +6.10%  '
$ws.Range('E49').Value = '  +8.03%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.09'
$ws.Range('E50').Value = '  +3.57%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.22'
$ws.Range('E51').Value = '  +11.60%  '
